$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 6786.3076
$ws.Range("I18").Value = 3929.2727
$ws.Range("K18").Value = 3929.2727
$ws.Range("M18").Value = -3645.2727
$ws.Range("H19").Value = 1653.9412
$ws.Range("I19").Value = 1559.5
$ws.Range("J19").Value = 1788.8572
$ws.Range("K19").Value = 1559.5
$ws.Range("L19").Value = 1788.8572
$ws.Range("M19").Value = -1384.5
$ws.Range("N19").Value = -2138.8572
$ws.Range("H107").Value = 1911.1428
$ws.Range("I107").Value = 1529.6666
$ws.Range("K107").Value = 1529.6666
$ws.Range("M107").Value = 390.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1912.9166
$ws.Range("I2").Value = 1782.5217
$ws.Range("K2").Value = 1782.5217
$ws.Range("M2").Value = -1669.5217
$ws.Range("H116").Value = 1912.9166
$ws.Range("I116").Value = 1782.5217
$ws.Range("K116").Value = 1782.5217
$ws.Range("M116").Value = 511.4783
$ws.Range("H127").Value = 96000
$ws.Range("J127").Value = 96000
$ws.Range("L127").Value = 96000
$ws.Range("N127").Value = -105920

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1912.9166
$ws.Range("I3").Value = 1782.5217
$ws.Range("K3").Value = 1782.5217
$ws.Range("M3").Value = -1668.5217
$ws.Range("H86").Value = 1443.2609
$ws.Range("I86").Value = 1206.4445
$ws.Range("J86").Value = 2295.8
$ws.Range("K86").Value = 1206.4445
$ws.Range("L86").Value = 2295.8
$ws.Range("M86").Value = -83.44450000000006
$ws.Range("N86").Value = -4541.8
$ws.Range("H89").Value = 1443.2609
$ws.Range("I89").Value = 1206.4445
$ws.Range("J89").Value = 2295.8
$ws.Range("K89").Value = 6032.2225
$ws.Range("L89").Value = 11479
$ws.Range("M89").Value = -416.2224999999999
$ws.Range("N89").Value = -22711
$ws.Range("H94").Value = 4313.846
$ws.Range("I94").Value = 6270.8125
$ws.Range("J94").Value = 1182.7
$ws.Range("K94").Value = 6270.8125
$ws.Range("L94").Value = 1182.7
$ws.Range("M94").Value = -5819.8125
$ws.Range("N94").Value = -2084.7
$ws.Range("H122").Value = 120000
$ws.Range("J122").Value = 120000
$ws.Range("L122").Value = 120000
$ws.Range("N122").Value = -129800
$ws.Range("H127").Value = 56823.715
$ws.Range("J127").Value = 56823.715
$ws.Range("L127").Value = 56823.715
$ws.Range("N127").Value = -66743.715
$ws.Range("H131").Value = 94880
$ws.Range("J131").Value = 94880
$ws.Range("L131").Value = 94880
$ws.Range("N131").Value = -104960

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20000794
$ws.Range("I31").Value = 20000794
$ws.Range("K31").Value = 20000794
$ws.Range("M31").Value = -20000499
$ws.Range("H34").Value = 20000794
$ws.Range("I34").Value = 20000794
$ws.Range("K34").Value = 20000794
$ws.Range("M34").Value = -20000592
$ws.Range("H58").Value = 1751.0834
$ws.Range("I58").Value = 1756.2667
$ws.Range("J58").Value = 1742.4445
$ws.Range("K58").Value = 1756.2667
$ws.Range("L58").Value = 1742.4445
$ws.Range("M58").Value = -1553.2667
$ws.Range("N58").Value = -2148.4445
$ws.Range("H100").Value = 38998.5
$ws.Range("J100").Value = 38998.5
$ws.Range("L100").Value = 38998.5
$ws.Range("N100").Value = -41162.5
$ws.Range("H105").Value = 1535.7693
$ws.Range("I105").Value = 906
$ws.Range("K105").Value = 906
$ws.Range("M105").Value = 841
$ws.Range("H125").Value = 95999.5
$ws.Range("J125").Value = 95999.5
$ws.Range("L125").Value = 95999.5
$ws.Range("N125").Value = -100919.5
$ws.Range("H132").Value = 78094.38
$ws.Range("I132").Value = 84477.25
$ws.Range("K132").Value = 253431.75
$ws.Range("M132").Value = -250901.75
$ws.Range("H134").Value = 2319.64
$ws.Range("I134").Value = 1633.2632
$ws.Range("K134").Value = 4899.7896
$ws.Range("M134").Value = -2364.7896
$ws.Range("H136").Value = 1751.0834
$ws.Range("I136").Value = 1756.2667
$ws.Range("J136").Value = 1742.4445
$ws.Range("K136").Value = 5268.800099999999
$ws.Range("L136").Value = 5227.333500000001
$ws.Range("M136").Value = -2718.800099999999
$ws.Range("N136").Value = -10327.3335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 155
$ws.Range("I14").Value = 155
$ws.Range("K14").Value = 465
$ws.Range("M14").Value = -292
$ws.Range("H131").Value = 612478.4399999999
$ws.Range("I131").Value = 612478.4399999999
$ws.Range("K131").Value = 1837435.32
$ws.Range("M131").Value = -1832395.32

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7599
$ws.Range("I70").Value = 7332
$ws.Range("K70").Value = 7332
$ws.Range("M70").Value = -7062
$ws.Range("H73").Value = 7599
$ws.Range("I73").Value = 7332
$ws.Range("K73").Value = 7332
$ws.Range("M73").Value = -6396
$ws.Range("H80").Value = 6307.375
$ws.Range("J80").Value = 10701.125
$ws.Range("L80").Value = 10701.125
$ws.Range("N80").Value = -12697.125
$ws.Range("H83").Value = 6307.375
$ws.Range("J83").Value = 10701.125
$ws.Range("L83").Value = 53505.625
$ws.Range("N83").Value = -63489.625
$ws.Range("H97").Value = 1438.5
$ws.Range("I97").Value = 1435.5834
$ws.Range("K97").Value = 1435.5834
$ws.Range("M97").Value = -939.5834
$ws.Range("H128").Value = 45000
$ws.Range("J128").Value = 45000
$ws.Range("L128").Value = 45000
$ws.Range("N128").Value = -54960

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4572.143
$ws.Range("I7").Value = 5200
$ws.Range("K7").Value = 5200
$ws.Range("M7").Value = -5088
$ws.Range("H38").Value = 89999.5
$ws.Range("J38").Value = 89999.5
$ws.Range("L38").Value = 89999.5
$ws.Range("N38").Value = -90819.5
$ws.Range("H46").Value = 3077.1052
$ws.Range("J46").Value = 4232.727
$ws.Range("L46").Value = 4232.727
$ws.Range("N46").Value = -4608.727
$ws.Range("H55").Value = 2597.75
$ws.Range("J55").Value = 3249.8333
$ws.Range("L55").Value = 3249.8333
$ws.Range("N55").Value = -3595.8333
$ws.Range("H68").Value = 3546
$ws.Range("I68").Value = 2932.5
$ws.Range("J68").Value = 6000
$ws.Range("K68").Value = 2932.5
$ws.Range("L68").Value = 6000
$ws.Range("M68").Value = -2183.5
$ws.Range("N68").Value = -7498
$ws.Range("H71").Value = 3546
$ws.Range("I71").Value = 2932.5
$ws.Range("J71").Value = 6000
$ws.Range("K71").Value = 14662.5
$ws.Range("L71").Value = 30000
$ws.Range("M71").Value = -10918.5
$ws.Range("N71").Value = -37488
$ws.Range("H82").Value = 9397.4
$ws.Range("I82").Value = 2181.3333
$ws.Range("J82").Value = 20221.5
$ws.Range("K82").Value = 2181.3333
$ws.Range("L82").Value = 20221.5
$ws.Range("M82").Value = -1820.3333
$ws.Range("N82").Value = -20943.5
$ws.Range("H85").Value = 9397.4
$ws.Range("I85").Value = 2181.3333
$ws.Range("J85").Value = 20221.5
$ws.Range("K85").Value = 2181.3333
$ws.Range("L85").Value = 20221.5
$ws.Range("M85").Value = -933.3332999999998
$ws.Range("N85").Value = -22717.5
$ws.Range("H126").Value = 4572.143
$ws.Range("I126").Value = 5200
$ws.Range("K126").Value = 15600
$ws.Range("M126").Value = -13130
$ws.Range("H136").Value = 4244.56
$ws.Range("I136").Value = 3743.4211
$ws.Range("K136").Value = 11230.2633
$ws.Range("M136").Value = -8680.263300000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3800.5
$ws.Range("I62").Value = 2851
$ws.Range("K62").Value = 2851
$ws.Range("M62").Value = -2227
$ws.Range("H65").Value = 3800.5
$ws.Range("I65").Value = 2851
$ws.Range("K65").Value = 14255
$ws.Range("M65").Value = -11135
$ws.Range("H136").Value = 17575
$ws.Range("I136").Value = 20485.47
$ws.Range("J136").Value = 3438.4285
$ws.Range("K136").Value = 61456.41
$ws.Range("L136").Value = 10315.2855
$ws.Range("M136").Value = -58906.41
$ws.Range("N136").Value = -15415.2855

Write-Output "Updated 205 cells across 8 sheets."